$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shorten all School District names in column A (rows 2-247)
$ws.Cells.Item(2, 1).Value = "Aberdeen"
$ws.Cells.Item(3, 1).Value = "Adna"
$ws.Cells.Item(4, 1).Value = "Anacortes"
$ws.Cells.Item(5, 1).Value = "Arlington"
$ws.Cells.Item(6, 1).Value = "Asotin-Anatone"
$ws.Cells.Item(7, 1).Value = "Auburn"
$ws.Cells.Item(8, 1).Value = "Bainbridge"
$ws.Cells.Item(9, 1).Value = "Battle Ground"
$ws.Cells.Item(10, 1).Value = "Bellevue"
$ws.Cells.Item(11, 1).Value = "Bellingham"
$ws.Cells.Item(12, 1).Value = "Bethel"
$ws.Cells.Item(13, 1).Value = "Bickleton"
$ws.Cells.Item(14, 1).Value = "Blaine"
$ws.Cells.Item(15, 1).Value = "Bremerton"
$ws.Cells.Item(16, 1).Value = "Brewster"
$ws.Cells.Item(17, 1).Value = "Bridgeport"
$ws.Cells.Item(18, 1).Value = "Burlington-Edison"
$ws.Cells.Item(19, 1).Value = "Camas"
$ws.Cells.Item(20, 1).Value = "Cape Flattery"
$ws.Cells.Item(21, 1).Value = "Cascade"
$ws.Cells.Item(22, 1).Value = "Cashmere"
$ws.Cells.Item(23, 1).Value = "Castle Rock"
$ws.Cells.Item(24, 1).Value = "Central Kitsap"
$ws.Cells.Item(25, 1).Value = "Central Valley"
$ws.Cells.Item(26, 1).Value = "Centralia"
$ws.Cells.Item(27, 1).Value = "Chehalis"
$ws.Cells.Item(28, 1).Value = "Cheney"
$ws.Cells.Item(29, 1).Value = "Chewelah"
$ws.Cells.Item(30, 1).Value = "Chimacum"
$ws.Cells.Item(31, 1).Value = "Clarkston"
$ws.Cells.Item(32, 1).Value = "Cle Elum-Roslyn"
$ws.Cells.Item(33, 1).Value = "Clover Park"
$ws.Cells.Item(34, 1).Value = "Colfax"
$ws.Cells.Item(35, 1).Value = "College Place"
$ws.Cells.Item(36, 1).Value = "Colton"
$ws.Cells.Item(37, 1).Value = "Columbia (Ste)"
$ws.Cells.Item(38, 1).Value = "Columbia (Wal)"
$ws.Cells.Item(39, 1).Value = "Colville"
$ws.Cells.Item(40, 1).Value = "Concrete"
$ws.Cells.Item(41, 1).Value = "Coulee-Hartline"
$ws.Cells.Item(42, 1).Value = "Coupeville"
$ws.Cells.Item(43, 1).Value = "Crescent"
$ws.Cells.Item(44, 1).Value = "Curlew"
$ws.Cells.Item(45, 1).Value = "Cusick"
$ws.Cells.Item(46, 1).Value = "Darrington"
$ws.Cells.Item(47, 1).Value = "Davenport"
$ws.Cells.Item(48, 1).Value = "Dayton"
$ws.Cells.Item(49, 1).Value = "Deer Park"
$ws.Cells.Item(50, 1).Value = "East Valley (Spo)"
$ws.Cells.Item(51, 1).Value = "East Valley (Yak)"
$ws.Cells.Item(52, 1).Value = "Eastmont"
$ws.Cells.Item(53, 1).Value = "Easton"
$ws.Cells.Item(54, 1).Value = "Eatonville"
$ws.Cells.Item(55, 1).Value = "Edmonds"
$ws.Cells.Item(56, 1).Value = "Ellensburg"
$ws.Cells.Item(57, 1).Value = "Elma"
$ws.Cells.Item(58, 1).Value = "Entiat"
$ws.Cells.Item(59, 1).Value = "Enumclaw"
$ws.Cells.Item(60, 1).Value = "Ephrata"
$ws.Cells.Item(61, 1).Value = "Everett"
$ws.Cells.Item(62, 1).Value = "Evergreen (Clark)"
$ws.Cells.Item(63, 1).Value = "Federal Way"
$ws.Cells.Item(64, 1).Value = "Ferndale"
$ws.Cells.Item(65, 1).Value = "Fife"
$ws.Cells.Item(66, 1).Value = "Finley"
$ws.Cells.Item(67, 1).Value = "Franklin Pierce"
$ws.Cells.Item(68, 1).Value = "Freeman"
$ws.Cells.Item(69, 1).Value = "Glenwood"
$ws.Cells.Item(70, 1).Value = "Goldendale"
$ws.Cells.Item(71, 1).Value = "Grand Coulee Dam"
$ws.Cells.Item(72, 1).Value = "Grandview"
$ws.Cells.Item(73, 1).Value = "Granger"
$ws.Cells.Item(74, 1).Value = "Granite Falls"
$ws.Cells.Item(75, 1).Value = "Harrington"
$ws.Cells.Item(76, 1).Value = "Highland"
$ws.Cells.Item(77, 1).Value = "Highline"
$ws.Cells.Item(78, 1).Value = "Hockinson"
$ws.Cells.Item(79, 1).Value = "Hoquiam"
$ws.Cells.Item(80, 1).Value = "Inchelium"
$ws.Cells.Item(81, 1).Value = "Issaquah"
$ws.Cells.Item(82, 1).Value = "Kahlotus"
$ws.Cells.Item(83, 1).Value = "Kalama"
$ws.Cells.Item(84, 1).Value = "Kelso"
$ws.Cells.Item(85, 1).Value = "Kennewick"
$ws.Cells.Item(86, 1).Value = "Kent"
$ws.Cells.Item(87, 1).Value = "Kettle Falls"
$ws.Cells.Item(88, 1).Value = "Kiona-Benton City"
$ws.Cells.Item(89, 1).Value = "Kittitas"
$ws.Cells.Item(90, 1).Value = "Klickitat"
$ws.Cells.Item(91, 1).Value = "La Center"
$ws.Cells.Item(92, 1).Value = "La Conner"
$ws.Cells.Item(93, 1).Value = "Lacrosse"
$ws.Cells.Item(94, 1).Value = "Lake Chelan"
$ws.Cells.Item(95, 1).Value = "Lake Quinault"
$ws.Cells.Item(96, 1).Value = "Lake Stevens"
$ws.Cells.Item(97, 1).Value = "Lake Washington"
$ws.Cells.Item(98, 1).Value = "Lakewood"
$ws.Cells.Item(99, 1).Value = "Liberty"
$ws.Cells.Item(100, 1).Value = "Longview"
$ws.Cells.Item(101, 1).Value = "Lopez Island"
$ws.Cells.Item(102, 1).Value = "Lyle"
$ws.Cells.Item(103, 1).Value = "Lynden"
$ws.Cells.Item(104, 1).Value = "Mabton"
$ws.Cells.Item(105, 1).Value = "Mansfield"
$ws.Cells.Item(106, 1).Value = "Manson"
$ws.Cells.Item(107, 1).Value = "Mary M Knight"
$ws.Cells.Item(108, 1).Value = "Mary Walker"
$ws.Cells.Item(109, 1).Value = "Marysville"
$ws.Cells.Item(110, 1).Value = "Mead"
$ws.Cells.Item(111, 1).Value = "Medical Lake"
$ws.Cells.Item(112, 1).Value = "Mercer Island"
$ws.Cells.Item(113, 1).Value = "Meridian"
$ws.Cells.Item(114, 1).Value = "Methow Valley"
$ws.Cells.Item(115, 1).Value = "Monroe"
$ws.Cells.Item(116, 1).Value = "Montesano"
$ws.Cells.Item(117, 1).Value = "Morton"
$ws.Cells.Item(118, 1).Value = "Moses Lake"
$ws.Cells.Item(119, 1).Value = "Mossyrock"
$ws.Cells.Item(120, 1).Value = "Mount Vernon"
$ws.Cells.Item(121, 1).Value = "Mt Adams"
$ws.Cells.Item(122, 1).Value = "Mt Baker"
$ws.Cells.Item(123, 1).Value = "Mukilteo"
$ws.Cells.Item(124, 1).Value = "Naches Valley"
$ws.Cells.Item(125, 1).Value = "Napavine"
$ws.Cells.Item(126, 1).Value = "Naselle-Grays R."
$ws.Cells.Item(127, 1).Value = "Newport"
$ws.Cells.Item(128, 1).Value = "Nine Mile Falls"
$ws.Cells.Item(129, 1).Value = "Nooksack Valley"
$ws.Cells.Item(130, 1).Value = "North Beach"
$ws.Cells.Item(131, 1).Value = "North Franklin"
$ws.Cells.Item(132, 1).Value = "North Kitsap"
$ws.Cells.Item(133, 1).Value = "North Mason"
$ws.Cells.Item(134, 1).Value = "North River"
$ws.Cells.Item(135, 1).Value = "North Thurston"
$ws.Cells.Item(136, 1).Value = "Northport"
$ws.Cells.Item(137, 1).Value = "Northshore"
$ws.Cells.Item(138, 1).Value = "Oak Harbor"
$ws.Cells.Item(139, 1).Value = "Oakesdale"
$ws.Cells.Item(140, 1).Value = "Ocean Beach"
$ws.Cells.Item(141, 1).Value = "Ocosta"
$ws.Cells.Item(142, 1).Value = "Odessa"
$ws.Cells.Item(143, 1).Value = "Okanogan"
$ws.Cells.Item(144, 1).Value = "Olympia"
$ws.Cells.Item(145, 1).Value = "Omak"
$ws.Cells.Item(146, 1).Value = "Onalaska"
$ws.Cells.Item(147, 1).Value = "Orcas Island"
$ws.Cells.Item(148, 1).Value = "Oroville"
$ws.Cells.Item(149, 1).Value = "Orting"
$ws.Cells.Item(150, 1).Value = "Othello"
$ws.Cells.Item(151, 1).Value = "Palouse"
$ws.Cells.Item(152, 1).Value = "Pasco"
$ws.Cells.Item(153, 1).Value = "Pateros"
$ws.Cells.Item(154, 1).Value = "Pe Ell"
$ws.Cells.Item(155, 1).Value = "Peninsula"
$ws.Cells.Item(156, 1).Value = "Pomeroy"
$ws.Cells.Item(157, 1).Value = "Port Angeles"
$ws.Cells.Item(158, 1).Value = "Port Townsend"
$ws.Cells.Item(159, 1).Value = "Prescott"
$ws.Cells.Item(160, 1).Value = "Prosser"
$ws.Cells.Item(161, 1).Value = "Pullman"
$ws.Cells.Item(162, 1).Value = "Puyallup"
$ws.Cells.Item(163, 1).Value = "Quilcene"
$ws.Cells.Item(164, 1).Value = "Quillayute Valley"
$ws.Cells.Item(165, 1).Value = "Quincy"
$ws.Cells.Item(166, 1).Value = "Rainier"
$ws.Cells.Item(167, 1).Value = "Raymond"
$ws.Cells.Item(168, 1).Value = "Reardan-Edwall"
$ws.Cells.Item(169, 1).Value = "Renton"
$ws.Cells.Item(170, 1).Value = "Republic"
$ws.Cells.Item(171, 1).Value = "Richland"
$ws.Cells.Item(172, 1).Value = "Ridgefield"
$ws.Cells.Item(173, 1).Value = "Ritzville"
$ws.Cells.Item(174, 1).Value = "Riverside"
$ws.Cells.Item(175, 1).Value = "Riverview"
$ws.Cells.Item(176, 1).Value = "Rochester"
$ws.Cells.Item(177, 1).Value = "Rosalia"
$ws.Cells.Item(178, 1).Value = "Royal"
$ws.Cells.Item(179, 1).Value = "Saint John"
$ws.Cells.Item(180, 1).Value = "San Juan Island"
$ws.Cells.Item(181, 1).Value = "Seattle"
$ws.Cells.Item(182, 1).Value = "Sedro Woolley"
$ws.Cells.Item(183, 1).Value = "Selah"
$ws.Cells.Item(184, 1).Value = "Selkirk"
$ws.Cells.Item(185, 1).Value = "Sequim"
$ws.Cells.Item(186, 1).Value = "Shelton"
$ws.Cells.Item(187, 1).Value = "Shoreline"
$ws.Cells.Item(188, 1).Value = "Skykomish"
$ws.Cells.Item(189, 1).Value = "Snohomish"
$ws.Cells.Item(190, 1).Value = "Snoqualmie Valley"
$ws.Cells.Item(191, 1).Value = "Soap Lake"
$ws.Cells.Item(192, 1).Value = "South Bend"
$ws.Cells.Item(193, 1).Value = "South Kitsap"
$ws.Cells.Item(194, 1).Value = "South Whidbey"
$ws.Cells.Item(195, 1).Value = "Spokane"
$ws.Cells.Item(196, 1).Value = "Sprague"
$ws.Cells.Item(197, 1).Value = "Stanwood-Camano"
$ws.Cells.Item(198, 1).Value = "Steilacoom Hist."
$ws.Cells.Item(199, 1).Value = "Stevenson-Carson"
$ws.Cells.Item(200, 1).Value = "Sultan"
$ws.Cells.Item(201, 1).Value = "Sumner"
$ws.Cells.Item(202, 1).Value = "Sunnyside"
$ws.Cells.Item(203, 1).Value = "Tacoma"
$ws.Cells.Item(204, 1).Value = "Taholah"
$ws.Cells.Item(205, 1).Value = "Tahoma"
$ws.Cells.Item(206, 1).Value = "Tekoa"
$ws.Cells.Item(207, 1).Value = "Tenino"
$ws.Cells.Item(208, 1).Value = "Thorp"
$ws.Cells.Item(209, 1).Value = "Toledo"
$ws.Cells.Item(210, 1).Value = "Tonasket"
$ws.Cells.Item(211, 1).Value = "Toppenish"
$ws.Cells.Item(212, 1).Value = "Touchet"
$ws.Cells.Item(213, 1).Value = "Toutle Lake"
$ws.Cells.Item(214, 1).Value = "Trout Lake"
$ws.Cells.Item(215, 1).Value = "Tukwila"
$ws.Cells.Item(216, 1).Value = "Tumwater"
$ws.Cells.Item(218, 1).Value = "University Place"
$ws.Cells.Item(219, 1).Value = "Valley"
$ws.Cells.Item(220, 1).Value = "Vancouver"
$ws.Cells.Item(221, 1).Value = "Vashon Island"
$ws.Cells.Item(222, 1).Value = "Wahkiakum"
$ws.Cells.Item(223, 1).Value = "Wahluke"
$ws.Cells.Item(224, 1).Value = "Waitsburg"
$ws.Cells.Item(225, 1).Value = "Walla Walla"
$ws.Cells.Item(226, 1).Value = "Wapato"
$ws.Cells.Item(227, 1).Value = "Warden"
$ws.Cells.Item(228, 1).Value = "Washougal"
$ws.Cells.Item(229, 1).Value = "Washtucna"
$ws.Cells.Item(230, 1).Value = "Waterville"
$ws.Cells.Item(231, 1).Value = "Wellpinit"
$ws.Cells.Item(232, 1).Value = "Wenatchee"
$ws.Cells.Item(233, 1).Value = "West Valley (Spo)"
$ws.Cells.Item(234, 1).Value = "West Valley (Yak)"
$ws.Cells.Item(235, 1).Value = "White Pass"
$ws.Cells.Item(236, 1).Value = "White River"
$ws.Cells.Item(237, 1).Value = "White Salmon Valley"
$ws.Cells.Item(238, 1).Value = "Wilbur"
$ws.Cells.Item(239, 1).Value = "Willapa Valley"
$ws.Cells.Item(240, 1).Value = "Wilson Creek"
$ws.Cells.Item(241, 1).Value = "Winlock"
$ws.Cells.Item(242, 1).Value = "Wishkah Valley"
$ws.Cells.Item(243, 1).Value = "Wishram"
$ws.Cells.Item(244, 1).Value = "Woodland"
$ws.Cells.Item(245, 1).Value = "Yakima"
$ws.Cells.Item(246, 1).Value = "Yelm Community"
$ws.Cells.Item(247, 1).Value = "Zillah"

# 2) Restyle the "plain" district-name cells: drop explicit Calibri font name
#    (falls back to workbook default font) and normalize wrap text.
$plainRange = $ws.Range("A2:A6,A9:A36,A39:A49,A52:A61,A63:A125,A127:A197,A200:A202,A204:A216,A219,A221:A232,A235:A247")
$plainRange.Font.Name = ""
$plainRange.WrapText = $false

# 3) A few district-name cells instead pick up the header's style (Auburn, Bainbridge,
#    Columbia x2, East Valley x2, Evergreen, Naselle-Grays, Steilacoom, Stevenson-Carson,
#    Tacoma, University Place, Vancouver, West Valley x2) via copy/paste of formats.
$ws.Range("A1").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
$ws.Range("A37:A38").PasteSpecial(-4122)
$ws.Range("A50:A51").PasteSpecial(-4122)
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A126").PasteSpecial(-4122)
$ws.Range("A198:A199").PasteSpecial(-4122)
$ws.Range("A203").PasteSpecial(-4122)
$ws.Range("A218").PasteSpecial(-4122)
$ws.Range("A220").PasteSpecial(-4122)
$ws.Range("A233:A234").PasteSpecial(-4122)
$excel.CutCopyMode = 0
